$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so values such as
# "1.00", "0.700" or "93.806.44" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "93.806.44"
$ws.Range("E2").Value = "  +4.39%  "

# Row 3
$ws.Range("D3").Value = "3.135.75"
$ws.Range("E3").Value = "  +0.70%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.48%  "

# Row 5
$ws.Range("D5").Value = "243.64"
$ws.Range("E5").Value = "  +3.65%  "

# Row 6
$ws.Range("D6").Value = "618.76"
$ws.Range("E6").Value = "  +0.74%  "

# Row 7
$ws.Range("D7").Value = "1.11"
$ws.Range("E7").Value = "  +1.94%  "

# Row 8
$ws.Range("E8").Value = "  +12.78%  "

# Row 9
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.15%  "

# Row 10
$ws.Range("D10").Value = "3.130.94"
$ws.Range("E10").Value = "  +0.73%  "

# Row 11
$ws.Range("D11").Value = "0.748"
$ws.Range("E11").Value = "  +3.03%  "

# Row 12
$ws.Range("E12").Value = "  +0.41%  "

# Row 13
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  +6.06%  "

# Row 14
$ws.Range("D14").Value = "35.02"
$ws.Range("E14").Value = "  +1.43%  "

# Row 15
$ws.Range("D15").Value = "93.179.70"
$ws.Range("E15").Value = "  +3.39%  "

# Row 16
$ws.Range("E16").Value = "  +0.79%  "

# Row 17
$ws.Range("D17").Value = "3.718.00"
$ws.Range("E17").Value = "  +0.51%  "

# Row 18
$ws.Range("D18").Value = "3.100.03"
$ws.Range("E18").Value = "  -1.71%  "

# Row 19
$ws.Range("D19").Value = "3.79"
$ws.Range("E19").Value = "  +4.34%  "

# Row 20
$ws.Range("D20").Value = "14.98"
$ws.Range("E20").Value = "  +0.93%  "

# Row 21
$ws.Range("E21").Value = "  +5.70%  "

# Row 22
$ws.Range("D22").Value = "5.88"
$ws.Range("E22").Value = "  +2.73%  "

# Row 23
$ws.Range("D23").Value = "9.46"
$ws.Range("E23").Value = "  +6.43%  "

# Row 24
$ws.Range("D24").Value = "452.60"
$ws.Range("E24").Value = "  +4.44%  "

# Row 25
$ws.Range("D25").Value = "5.91"
$ws.Range("E25").Value = "  +4.69%  "

# Row 26
$ws.Range("D26").Value = "88.13"
$ws.Range("E26").Value = "  +8.06%  "

# Row 27
$ws.Range("E27").Value = "  +2.81%  "

# Row 28
$ws.Range("D28").Value = "3.294.87"
$ws.Range("E28").Value = "  -1.46%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("D30").Value = "0.138"
$ws.Range("E30").Value = "  +10.38%  "

# Row 31
$ws.Range("D31").Value = "0.171"
$ws.Range("E31").Value = "  +1.79%  "

# Row 32
$ws.Range("E32").Value = "  +0.50%  "

# Row 33
$ws.Range("D33").Value = "9.32"
$ws.Range("E33").Value = "  +2.01%  "

# Row 34
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D34").Value = "8.27"
$ws.Range("E34").Value = "  +9.12%  "

# Row 35
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  +9.55%  "

# Row 36
$ws.Range("D36").Value = "0.162"
$ws.Range("E36").Value = "  -3.41%  "

# Row 37
$ws.Range("D37").Value = "26.53"
$ws.Range("E37").Value = "  +3.16%  "

# Row 38
$ws.Range("D38").Value = "1.93"
$ws.Range("E38").Value = "  +0.96%  "

# Row 39
$ws.Range("D39").Value = "3.91"
$ws.Range("E39").Value = "  +4.60%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "485.16"
$ws.Range("E40").Value = "  -1.26%  "

# Row 41
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "1.32"
$ws.Range("E41").Value = "  +0.07%  "

# Row 42
$ws.Range("E42").Value = "  -0.09%  "

# Row 43
$ws.Range("E43").Value = "  +2.55%  "

# Row 44
$ws.Range("D44").Value = "23.10"
$ws.Range("E44").Value = "  +4.71%  "

# Row 46
$ws.Range("D46").Value = "162.26"
$ws.Range("E46").Value = "  +2.14%  "

# Row 47
$ws.Range("E47").Value = "  +4.28%  "

# Row 48
$ws.Range("D48").Value = "0.700"
$ws.Range("E48").Value = "  +0.37%  "

# Row 49
$ws.Range("E49").Value = "  +5.35%  "

# Row 50
$ws.Range("D50").Value = "0.0336"
$ws.Range("E50").Value = "  +6.82%  "

# Row 51
$ws.Range("D51").Value = "4.52"
$ws.Range("E51").Value = "  +4.25%  "
